# Adds a new "Dark_light" parameter column (F) to the metadata sheet,
# classifying every sample row as light ("l") or dark ("d") based on its
# Timepoints value (column C), and tweaks a couple of view/format bits
# that tagged along with the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header -----------------------------------------------
$ws.Range("F1").Value = "Dark_light"

# --- Rows sampled during the "light" phase ----------------------------
$lightRows = @(
    2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,
    38,39,40,41,42,43,44,45,46,47,48,49,50,51,
    64,65,66,67,68,69,70,71,72,73,74,75
)

# --- Rows sampled during the "dark" phase -----------------------------
$darkRows = @(
    19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,
    52,53,54,55,56,57,58,59,60,61,62,63
)

foreach ($r in $lightRows) {
    $ws.Cells.Item($r, 6).Value = "l"
}

foreach ($r in $darkRows) {
    $ws.Cells.Item($r, 6).Value = "d"
}

# --- Incidental time-format touches on a few empty H cells -------------
$ws.Range("H2").NumberFormat = "h:mm"
$ws.Range("H3:H10").NumberFormat = "[$-F400]h:mm:ss\ AM/PM"

# --- Selection left where the author's cursor ended up ------------------
$ws.Range("J9").Select()
